# Added a new case - AddToCart, and added common actions.
#
# - Select A3 and zoom to 160% on the existing "ableToSearchProducts" sheet
#   (it becomes the non-active tab once the new sheet below is added).
# - Add a new worksheet "ableToAddToCart" right after it, with Qty/Nikon
#   D300/iPod Shuffle/$122.00 data, and make it the active tab.

function Set-TextValue($cell, $text) {
    # Forces the value to be stored as a string (shared-string) cell even
    # when it looks numeric/currency (e.g. "$98.00"), then strips the
    # Text-number-format styling that gets attached along the way so the
    # cell ends up with the default style.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Common action: select A3 on the first sheet and bump the zoom level.
$ws1.Range("A3").Select()
$excel.ActiveWindow.Zoom = 160

# New case: AddToCart sheet, inserted right after the search-products sheet.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "ableToAddToCart"

$ws2.Range("A1").Value = "Product Name"
$ws2.Range("B1").Value = "Product  Price"
$ws2.Range("C1").Value = "Qty"

$ws2.Range("A2").Value = "Nikon D300"
Set-TextValue $ws2.Range("B2") "`$98.00"
$ws2.Range("C2").Value = 3

$ws2.Range("A3").Value = "iPod Shuffle"
Set-TextValue $ws2.Range("B3") "`$122.00"
$ws2.Range("C3").Value = 2

$ws2.Range("A4").Value = "Samsung Galaxy Tab 10.1"
Set-TextValue $ws2.Range("B4") "`$241.99"
$ws2.Range("C4").Value = 1
